$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 121210411.712973
$ws.Range("D2").Value = 57.495063

$ws.Range("B3").Value = 20378896.423127
$ws.Range("D3").Value = 4.833273
$ws.Range("E3").Value = 0.008814000000000001

$ws.Range("B4").Value = 470125964.657988
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = -567.933683
$ws.Range("H5").Value = -1202.057117
$ws.Range("I5").Value = 66.189751
$ws.Range("J5").Value = 0.08945400000000001

$ws.Range("G6").Value = 49.487735
$ws.Range("H6").Value = -633.328937
$ws.Range("I6").Value = 732.304408
$ws.Range("J6").Value = 0.984009

$ws.Range("G7").Value = 617.421418
$ws.Range("H7").Value = 101.204913
$ws.Range("I7").Value = 1133.637924
$ws.Range("J7").Value = 0.014345
